$wb = $excel.ActiveWorkbook

# --- Add the new sheet as a copy of "Putz 204" (the closest matching layout:
# same column widths / borders / merges), placed at the very end of the
# workbook, then rename + re-populate it for the "El Shaddai 301" delivery
# challan. Copying (instead of Worksheets.Add()) preserves the exact column
# widths/row heights/merge cells/borders that a brand-new sheet would not
# get automatically. ---
$src = $wb.Worksheets.Item("Putz 204")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "El Shaddai 301"

# --- Header block ---
$ws.Range("A1").Value = "DC  no:23-24QEl Shaddai 301"
$ws.Range("E1").Value = "date: 26/12/2023"

$ws.Range("A2").Value = "client:  El Shaddai"
$ws.Range("E2").Value = "Ref: Delivery of Material                                                                  "

$ws.Range("A3").Value = "Billing Address:  El Shaddai"
$ws.Range("E3").Value = ""

$ws.Range("A5").Value = "Shipping Address: El Shaddai Mapusa Goa"

# --- Table header row (row 7): bold 10pt Calibri, add a QTY column ---
$ws.Range("B7").Value = "ITEM DESCRIPTION"
$ws.Range("C7").Value = "QTY"
$headerRng = $ws.Range("B7:C7")
$headerRng.Font.Name = "Calibri  "
$headerRng.Font.Size = 10
$headerRng.Font.Bold = $true
$headerRng.Font.Color = $ws.Range("A7").Font.Color
$headerRng.HorizontalAlignment = $ws.Range("A7").HorizontalAlignment
$headerRng.VerticalAlignment = $ws.Range("A7").VerticalAlignment
$headerRng.WrapText = $true
$headerRng.Borders.LineStyle = $ws.Range("A7").Borders.LineStyle

# --- Table data rows 8-10: regular 10pt Calibri ---
$ws.Range("B8").Value = "Fuji 10 Kva UPS"
$ws.Range("C8").Value = 1

$ws.Range("B9").Value = "26 AH Battery"
$ws.Range("C9").Value = 16

# Insert a brand-new row 10 (pushes the trailing "terms" row from 11 to 12,
# matching the target layout) and fill it in.
$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Rack"
$ws.Range("C10").Value = 1
$ws.Rows.Item(10).RowHeight = $ws.Rows.Item(9).RowHeight

$dataRng = $ws.Range("B8:C10")
$dataRng.Font.Name = "Calibri  "
$dataRng.Font.Size = 10
$dataRng.Font.Bold = $false
$dataRng.Font.Color = $ws.Range("A8").Font.Color
$dataRng.HorizontalAlignment = $ws.Range("A8").HorizontalAlignment
$dataRng.VerticalAlignment = $ws.Range("A8").VerticalAlignment
$dataRng.WrapText = $true
$dataRng.Borders.LineStyle = $ws.Range("A8").Borders.LineStyle

$ws.Range("A10").HorizontalAlignment = $ws.Range("A9").HorizontalAlignment
$ws.Range("A10").VerticalAlignment = $ws.Range("A9").VerticalAlignment
$ws.Range("A10").Borders.LineStyle = $ws.Range("A9").Borders.LineStyle

# Leave row 11 blank as a spacer (matches the template) - row 12 keeps the
# "Terms and conditions..." note that shifted down when row 10 was inserted.

# The newly active sheet becomes the selected tab; the previous last sheet
# ("Putz 204") loses its tab-selection + goes back to a plain "select all" state.
$src.Cells.Select()
$ws.Range("A1").Select()

$wb.Save()
